$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.043.69'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '1.620.55'
$ws.Range("E3").Value = '  -1.58%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '214.71'
$ws.Range("E5").Value = '  -1.27%  '
$ws.Range("D6").Value = '0.518'
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '0.253'
$ws.Range("E8").Value = '  -0.78%  '
$ws.Range("D9").Value = '0.0628'
$ws.Range("E9").Value = '  -0.30%  '
$ws.Range("E10").Value = '  -0.16%  '
$ws.Range("D11").Value = '0.0849'
$ws.Range("E11").Value = '  -0.22%  '
$ws.Range("D12").Value = '1.622.43'
$ws.Range("E12").Value = '  -1.66%  '
$ws.Range("D13").Value = '4.16'
$ws.Range("E13").Value = '  -0.02%  '
$ws.Range("E14").Value = '  -1.38%  '
$ws.Range("D15").Value = '64.78'
$ws.Range("E15").Value = '  -3.71%  '
$ws.Range("D16").Value = '27.036.57'
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("D17").Value = '0.0₃0751'
$ws.Range("E17").Value = '  +1.29%  '
$ws.Range("D18").Value = '214.69'
$ws.Range("E18").Value = '  -2.62%  '
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D20").Value = '6.88'
$ws.Range("E20").Value = '  -1.03%  '
$ws.Range("D21").Value = '4.38'
$ws.Range("E21").Value = '  -0.86%  '
$ws.Range("D22").Value = '2.39'
$ws.Range("E22").Value = '  -5.61%  '
$ws.Range("E23").Value = '  -1.13%  '
$ws.Range("D24").Value = '148.37'
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").Value = '7.44'
$ws.Range("E25").Value = '  -0.85%  '
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("E27").Value = '  -1.86%  '
$ws.Range("D28").Value = '15.59'
$ws.Range("E28").Value = '  -1.25%  '
$ws.Range("E29").Value = '  +0.91%  '
$ws.Range("E30").Value = '  -0.95%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '3.37'
$ws.Range("E31").Value = '  -0.76%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '0.758'
$ws.Range("E32").Value = '  +36.88%  '
$ws.Range("E33").Value = '  -0.46%  '
$ws.Range("D34").Value = '1.345.55'
$ws.Range("E34").Value = '  +3.61%  '
$ws.Range("D35").Value = '1.57'
$ws.Range("E35").Value = '  -0.84%  '
$ws.Range("E36").Value = '  -0.55%  '
$ws.Range("E37").Value = '  +1.04%  '
$ws.Range("D38").Value = '0.855'
$ws.Range("E38").Value = '  -1.02%  '
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("D40").Value = '0.802'
$ws.Range("E40").Value = '  -1.07%  '
$ws.Range("D41").Value = '65.36'
$ws.Range("E41").Value = '  +4.82%  '
$ws.Range("D42").Value = '2.23'
$ws.Range("E42").Value = '  -0.51%  '
$ws.Range("D43").Value = '5.33'
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").Value = '1.758.32'
$ws.Range("E44").Value = '  -1.70%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '90.02'
$ws.Range("E45").Value = '  -2.35%  '
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '0.880'
$ws.Range("E46").Value = '  +30.98%  '
$ws.Range("D47").Value = '1.64'
$ws.Range("E47").Value = '  +1.99%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0107'
$ws.Range("E48").Value = '  +0.88%  '
$ws.Range("D49").Value = '0.101'
$ws.Range("E49").Value = '  +3.61%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.0516'
$ws.Range("E50").Value = '  +0.46%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '7.66'
$ws.Range("E51").Value = '  -0.94%  '
